# Regenerate save_data column G ("K") values: switch from old Strike# derived
# values to the newly computed K values (recomputed std/mean, calc and write
# s_vals). This updates column G (K) for rows 2-52 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2 = 1
    3 = 3
    4 = 1
    5 = 1
    6 = 1
    7 = 0
    8 = 2
    9 = 2
    10 = 1
    11 = 1
    12 = 0
    13 = 0
    14 = 1
    15 = 0
    16 = 1
    17 = 0
    18 = 1
    19 = 1
    20 = 2
    21 = 2
    22 = 0
    23 = 1
    24 = 2
    25 = 0
    26 = 4
    27 = 2
    28 = 3
    29 = 1
    30 = 3
    31 = 2
    32 = 1
    33 = 1
    34 = 0
    35 = 1
    37 = 1
    39 = 1
    40 = 4
    42 = 4
    43 = 1
    44 = 0
    46 = 1
    47 = 1
    48 = 0
    49 = 1
    50 = 2
    51 = 2
    52 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
